# Apply "Update countries & provincias Spain" edit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the "last updated" timestamp banner in A1
$ws.Range("A1").Value = "Datos actualizados a 11 de Agosto de 2020 a las 16:48"

# Country rows with updated case counts. A few pairs of adjacent rows also swap which
# country occupies which row (matching the shared-strings reorder in the source diff):
#   row 99/100  : Albania now before Mauritania
#   row 117/118 : Namibia now before Somalia
#   row 176-178 : Trinidad yTobago now before Mongolia/Eritrea
#   row 213/214 : Islas Malvinas now before Montserrat

# Row 4: Estados Unidos
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 5254561
$ws.Range("C4").Value = 3115
$ws.Range("D4").Value = 2717513
$ws.Range("E4").Value = 2370753
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 103
$ws.Range("H4").Value = 166295

# Row 6: India
$ws.Range("A6").Value = "India"
$ws.Range("B6").Value = 2294438
$ws.Range("C6").Value = 27285
$ws.Range("D6").Value = 1604119
$ws.Range("E6").Value = 644722
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 244
$ws.Range("H6").Value = 45597

# Row 19: Argentina
$ws.Range("A19").Value = "Argentina"
$ws.Range("B19").Value = 253868
$ws.Range("C19").Value = 0
$ws.Range("D19").Value = 181389
$ws.Range("E19").Value = 67694
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 21
$ws.Range("H19").Value = 4785

# Row 22: Alemania
$ws.Range("A22").Value = "Alemania"
$ws.Range("B22").Value = 218852
$ws.Range("C22").Value = 352
$ws.Range("D22").Value = 198900
$ws.Range("E22").Value = 10686
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 9266

# Row 49: Portugal
$ws.Range("A49").Value = "Portugal"
$ws.Range("B49").Value = 52945
$ws.Range("C49").Value = 120
$ws.Range("D49").Value = 38760
$ws.Range("E49").Value = 12424
$ws.Range("F49").Value = 0
$ws.Range("G49").Value = 2
$ws.Range("H49").Value = 1761

# Row 64: Moldavia
$ws.Range("A64").Value = "Moldavia"
$ws.Range("B64").Value = 28223
$ws.Range("C64").Value = 382
$ws.Range("D64").Value = 19740
$ws.Range("E64").Value = 7626
$ws.Range("F64").Value = 0
$ws.Range("G64").Value = 7
$ws.Range("H64").Value = 857

# Row 86: Noruega
$ws.Range("A86").Value = "Noruega"
$ws.Range("B86").Value = 9712
$ws.Range("C86").Value = 28
$ws.Range("D86").Value = 8857
$ws.Range("E86").Value = 599
$ws.Range("F86").Value = 0
$ws.Range("G86").Value = 0
$ws.Range("H86").Value = 256

# Row 93: Tayikistan
$ws.Range("A93").Value = "Tayikistan"
$ws.Range("B93").Value = 7871
$ws.Range("C93").Value = 44
$ws.Range("D93").Value = 6653
$ws.Range("E93").Value = 1155
$ws.Range("F93").Value = 0
$ws.Range("G93").Value = 1
$ws.Range("H93").Value = 63

# Row 99: Albania
$ws.Range("A99").Value = "Albania"
$ws.Range("B99").Value = 6676
$ws.Range("C99").Value = 140
$ws.Range("D99").Value = 3480
$ws.Range("E99").Value = 2991
$ws.Range("F99").Value = 0
$ws.Range("G99").Value = 5
$ws.Range("H99").Value = 205

# Row 100: Mauritania
$ws.Range("A100").Value = "Mauritania"
$ws.Range("B100").Value = 6555
$ws.Range("C100").Value = 0
$ws.Range("D100").Value = 5570
$ws.Range("E100").Value = 828
$ws.Range("F100").Value = 0
$ws.Range("G100").Value = 0
$ws.Range("H100").Value = 157

# Row 117: Namibia
$ws.Range("A117").Value = "Namibia"
$ws.Range("B117").Value = 3229
$ws.Range("C117").Value = 128
$ws.Range("D117").Value = 715
$ws.Range("E117").Value = 2495
$ws.Range("F117").Value = 0
$ws.Range("G117").Value = 0
$ws.Range("H117").Value = 19

# Row 118: Somalia
$ws.Range("A118").Value = "Somalia"
$ws.Range("B118").Value = 3227
$ws.Range("C118").Value = 0
$ws.Range("D118").Value = 1728
$ws.Range("E118").Value = 1406
$ws.Range("F118").Value = 0
$ws.Range("G118").Value = 0
$ws.Range("H118").Value = 93

# Row 170: Birmania
$ws.Range("A170").Value = "Birmania"
$ws.Range("B170").Value = 360
$ws.Range("C170").Value = 0
$ws.Range("D170").Value = 313
$ws.Range("E170").Value = 41
$ws.Range("F170").Value = 0
$ws.Range("G170").Value = 0
$ws.Range("H170").Value = 6

# Row 176: Trinidad yTobago
$ws.Range("A176").Value = "Trinidad yTobago"
$ws.Range("B176").Value = 295
$ws.Range("C176").Value = 14
$ws.Range("D176").Value = 139
$ws.Range("E176").Value = 148
$ws.Range("F176").Value = 0
$ws.Range("G176").Value = 0
$ws.Range("H176").Value = 8

# Row 177: Mongolia
$ws.Range("A177").Value = "Mongolia"
$ws.Range("B177").Value = 293
$ws.Range("C177").Value = 0
$ws.Range("D177").Value = 263
$ws.Range("E177").Value = 30
$ws.Range("F177").Value = 0
$ws.Range("G177").Value = 0
$ws.Range("H177").Value = 0

# Row 178: Eritrea
$ws.Range("A178").Value = "Eritrea"
$ws.Range("B178").Value = 285
$ws.Range("C178").Value = 0
$ws.Range("D178").Value = 248
$ws.Range("E178").Value = 37
$ws.Range("F178").Value = 0
$ws.Range("G178").Value = 0
$ws.Range("H178").Value = 0

# Row 190: Seychelles
$ws.Range("A190").Value = "Seychelles"
$ws.Range("B190").Value = 127
$ws.Range("C190").Value = 1
$ws.Range("D190").Value = 126
$ws.Range("E190").Value = 1
$ws.Range("F190").Value = 0
$ws.Range("G190").Value = 0
$ws.Range("H190").Value = 0

# Row 194: Liechtenstein
$ws.Range("A194").Value = "Liechtenstein"
$ws.Range("B194").Value = 89
$ws.Range("C194").Value = 0
$ws.Range("D194").Value = 87
$ws.Range("E194").Value = 1
$ws.Range("F194").Value = 0
$ws.Range("G194").Value = 0
$ws.Range("H194").Value = 1

# Row 213: Islas Malvinas
$ws.Range("A213").Value = "Islas Malvinas"
$ws.Range("B213").Value = 13
$ws.Range("C213").Value = 0
$ws.Range("D213").Value = 13
$ws.Range("E213").Value = 0
$ws.Range("F213").Value = 0
$ws.Range("G213").Value = 0
$ws.Range("H213").Value = 0

# Row 214: Montserrat
$ws.Range("A214").Value = "Montserrat"
$ws.Range("B214").Value = 13
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 12
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 1
